$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 was duplicating Row 2's data (Ang, Bryan). Update it to reflect
# the correct student (Wang Ma, Frank), matching rows 4/5.
$ws.Range("A3").Value = "Wang Ma, Frank"
$ws.Range("F3").Value = "fwan175@aucklanduni.ac.nz"
$ws.Range("G3").Value = 184846458
